{"js": "// The document contains six `<id>...</id>` markers, each currently split\n// across three runs as `<id>` + `p053r_aN` + `</id>` (the middle run holds\n// the \"aN\" suffix). The edit collapses each of those three runs into a\n// single run whose text is the merged `<id>p053r_N</id>` string (dropping\n// the \"a\" in the id), keeping the formatting of the first (\"<id>\") run.\nconst body = context.document.body;\n\nfor (let n = 1; n <= 6; n++) {\n  const oldText = `<id>p053r_a${n}</id>`;\n  const newText = `<id>p053r_${n}</id>`;\n\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    continue;\n  }\n\n  // Replacing the whole matched range's text merges the three runs it\n  // spans into a single run, carrying over the formatting of the first\n  // run in the range (the Courier New / 7f6000 \"<id>\" run).\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document contains six `<id>...</id>` markers, each currently split\n# across three runs as `<id>` + `p053r_aN` + `</id>` (the middle run holds\n# the \"aN\" suffix). The edit collapses each of those three runs into a\n# single run whose text is the merged `<id>p053r_N</id>` string (dropping\n# the \"a\" in the id), keeping the formatting of the first (\"<id>\") run.\n#\n# Searching for, and replacing, the *entire* \"<id>p053r_aN</id>\" span (not\n# just the \"p053r_aN\" word) makes Find.Execute's replacement collapse the\n# three runs it covers into a single run - matching the target edit.\n$d = $word.ActiveDocument\n\nfor ($n = 1; $n -le 6; $n++) {\n    $oldText = \"<id>p053r_a$n</id>\"\n    $newText = \"<id>p053r_$n</id>\"\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n\n    # wdFindContinue = 1, wdReplaceOne = 1\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n}\n"}
